# Split the run "), GUI.changed = true; .." into:
#   - "), GUI.changed = true;"                         (trailing " .." removed)
#   - a new run with " Selection.selectionChanged …"    (appended right after)
#
# We locate the exact run via its distinctive text, then use Range.InsertXML
# on exactly that run's extent so that the surrounding (unrelated) runs in
# the paragraph are left completely untouched (a plain text/Find replace
# would cause the engine to coalesce neighbouring same-format runs).

$d = $word.ActiveDocument

$oldRunText = "), GUI.changed = true; .."
$newFirstText = "), GUI.changed = true;"
$newSecondText = " Selection.selectionChanged " + [char]0x2026

foreach ($p in $d.Paragraphs) {
    $paraText = $p.Range.Text
    if ($paraText.Contains($oldRunText)) {
        $paraStart = $p.Range.Start
        $offset = $paraText.IndexOf($oldRunText)

        $targetStart = $paraStart + $offset
        $targetEnd = $targetStart + $oldRunText.Length

        $target = $d.Range($targetStart, $targetEnd)

        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body>' +
               '<w:p>' +
               '<w:r w:rsidR="005F0CD3">' +
               '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
               '<w:t>' + $newFirstText + '</w:t>' +
               '</w:r>' +
               '<w:r>' +
               '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
               '<w:t xml:space="preserve">' + $newSecondText + '</w:t>' +
               '</w:r>' +
               '</w:p>' +
               '</w:body>' +
               '</w:document>' +
               '</pkg:xmlData></pkg:part></pkg:package>'

        $target.InsertXML($xml)
        break
    }
}
